$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.136.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.124.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.132.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.623.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "363.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("E22").Value = "  -3.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.60%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.59%  "

$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "555.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.78%  "

$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("E34").Value = "  -2.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("E40").Value = "  -3.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("E42").Value = "  +4.75%  "

$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("E44").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.97%  "

$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("E51").Value = "  -0.04%  "

